# Cookie checker relative path
#
# The checker now also finds "WSS_FullScreenMode" -- a cookie that is set on
# a relative path and therefore shows up against (almost) every URL in the
# sheet instead of just one. Each affected row gets its cookie COUNT (col B)
# bumped by one and "WSS_FullScreenMode" appended as the new last cookie
# name in the row (after whatever cookies were already listed there).
#
# Rows 9 ("https://www.act4greece.gr/"), 19 ("https://www.waiz.gr/") and 24
# ("https://www.i-fund.gr/public/#!/") are unaffected and keep their
# original values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that previously had no cookies at all -> now just WSS_FullScreenMode.
$noCookieRows = @(2,3,4,5,6)
foreach ($r in $noCookieRows) {
    $ws.Range("B$r").Value = 1
    $ws.Range("C$r").Value = "WSS_FullScreenMode"
}

# Rows that previously had only "NBGPublicSite" -> append WSS_FullScreenMode.
$publicSiteRows = @(7,8,10,11,12,13,14,15)
foreach ($r in $publicSiteRows) {
    $ws.Range("B$r").Value = 2
    $ws.Range("C$r").Value = "NBGPublicSite"
    $ws.Range("D$r").Value = "WSS_FullScreenMode"
}

# Rows that previously had "NBGPUBLICConsent","NBGPublicSite" -> append WSS_FullScreenMode.
$consentRows = @(16,17,18,20,21,22,23,25,26,27,28)
foreach ($r in $consentRows) {
    $ws.Range("B$r").Value = 3
    $ws.Range("C$r").Value = "NBGPUBLICConsent"
    $ws.Range("D$r").Value = "NBGPublicSite"
    $ws.Range("E$r").Value = "WSS_FullScreenMode"
}

# Rows that previously had "NBGPUBLICConsent","NBGPublicSite","_af" -> insert
# WSS_FullScreenMode before the trailing "_af" cookie.
$afRows = @(29,30,31,32,33,34)
foreach ($r in $afRows) {
    $ws.Range("B$r").Value = 4
    $ws.Range("C$r").Value = "NBGPUBLICConsent"
    $ws.Range("D$r").Value = "NBGPublicSite"
    $ws.Range("E$r").Value = "WSS_FullScreenMode"
    $ws.Range("F$r").Value = "_af"
}
